# Update database (shift yearly columns by one year, dropping the oldest
# "1396/12" period and appending a new "1401/12" period) and refresh the
# figures for the affected metrics per the new read_price algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Columns E..I hold, in order, the 5 most-recent twelve-month periods.
$cols = @(5, 6, 7, 8, 9)

# --- Header labels (row 8 and row 24): shift the "N months ended" labels
# one column to the left and append the new trailing period.
$headers = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(8, $cols[$i]).Value = $headers[$i]
    $ws.Cells.Item(24, $cols[$i]).Value = $headers[$i]
}

# --- Row 16: هزینه استهلاک (depreciation expense)
$row16 = @(1728, 8847, 1493, 1638, 1780)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(16, $cols[$i]).Value = $row16[$i]
}

# --- Row 17: هزینه حقوق و دستمزد (payroll expense)
$row17 = @(60470, 74129, 74129, 128692, 229045)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(17, $cols[$i]).Value = $row17[$i]
}

# --- Row 19: سایر هزینه ها (other expenses)
$row19 = @(136342, 51001, 64747, 174860, 108729)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(19, $cols[$i]).Value = $row19[$i]
}

# --- Row 20: جمع (total)
$row20 = @(198540, 133977, 140369, 305190, 339554)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(20, $cols[$i]).Value = $row20[$i]
}

# --- Row 26: تعداد پرسنل غیر تولیدی شرکت (non-production personnel count)
$row26 = @(277, 274, 271, 273, 272)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(26, $cols[$i]).Value = $row26[$i]
}

# --- Row 27: تعداد پرسنل تولیدی شرکت (production personnel count)
$row27 = @(212, 243, 258, 253, 243)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item(27, $cols[$i]).Value = $row27[$i]
}
